{"js": "// Apply the \"Added many more features\" edits to the Blood Moon Wilds review.\n// Each entry is an exact, case-sensitive search string paired with its\n// replacement. Office.js `search()` matches ranges within single runs, and\n// `insertText(text, \"Replace\")` swaps the matched range's text while keeping\n// the run's existing formatting (bold/italic/heading style, etc.) intact.\nconst replacements = [\n  // Title appears twice verbatim: once as the Heading1 at the top of the\n  // document, and once again as a bold run near the bottom.\n  [\n    \"Play Blood Moon Wilds Free Today - Exciting Werewolf Slot Game\",\n    \"Play Blood Moon Wilds Free - Review of this Exciting Slot Game\"\n  ],\n  // \"What we like\" bullet list.\n  [\n    \"Lunar Calendar feature can turn high-value symbols into wilds\",\n    \"Lunar Calendar feature adds excitement to gameplay\"\n  ],\n  [\n    \"Four random functions that can result in extra wilds, multipliers, scatters, and Lunar Calendar triggers\",\n    \"Random functions offer additional chances to win\"\n  ],\n  [\n    \"Well-designed graphics and spooky sound effects that add to the eerie atmosphere\",\n    \"Impressive graphics and design create a dark and eerie atmosphere\"\n  ],\n  [\n    \"Compatible with all devices, including desktop computers, smartphones, and tablets\",\n    \"Compatible with all devices for convenient play\"\n  ],\n  // \"What we don't like\" bullet list.\n  [\n    \"Limited betting range of \\u00a30.10-\\u00a3100 per spin\",\n    \"Limited number of paylines\"\n  ],\n  [\n    \"Free spins feature can be difficult to trigger\",\n    \"May not appeal to players who prefer lighter themes\"\n  ],\n  // Closing italic summary line.\n  [\n    \"Try Blood Moon Wilds slot game for free today and discover exciting werewolf characters, lunar calendar feature, and more. Compatible on all devices.\",\n    \"Play Blood Moon Wilds for free and explore its exciting features and eerie atmosphere.\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edits to the Blood Moon Wilds review\n# using Word COM interop (Find/Replace). Find.Execute with Replace=2\n# (wdReplaceAll) rewrites every matching run's text in one call while\n# preserving each run's existing character formatting (bold/italic/\n# heading style, etc.).\n\n$d = $word.ActiveDocument\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @{\n        Old = \"Play Blood Moon Wilds Free Today - Exciting Werewolf Slot Game\"\n        New = \"Play Blood Moon Wilds Free - Review of this Exciting Slot Game\"\n    },\n    @{\n        Old = \"Lunar Calendar feature can turn high-value symbols into wilds\"\n        New = \"Lunar Calendar feature adds excitement to gameplay\"\n    },\n    @{\n        Old = \"Four random functions that can result in extra wilds, multipliers, scatters, and Lunar Calendar triggers\"\n        New = \"Random functions offer additional chances to win\"\n    },\n    @{\n        Old = \"Well-designed graphics and spooky sound effects that add to the eerie atmosphere\"\n        New = \"Impressive graphics and design create a dark and eerie atmosphere\"\n    },\n    @{\n        Old = \"Compatible with all devices, including desktop computers, smartphones, and tablets\"\n        New = \"Compatible with all devices for convenient play\"\n    },\n    @{\n        Old = \"Limited betting range of \u00a30.10-\u00a3100 per spin\"\n        New = \"Limited number of paylines\"\n    },\n    @{\n        Old = \"Free spins feature can be difficult to trigger\"\n        New = \"May not appeal to players who prefer lighter themes\"\n    },\n    @{\n        Old = \"Try Blood Moon Wilds slot game for free today and discover exciting werewolf characters, lunar calendar feature, and more. Compatible on all devices.\"\n        New = \"Play Blood Moon Wilds for free and explore its exciting features and eerie atmosphere.\"\n    }\n)\n\nforeach ($rep in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($rep.Old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $rep.New, $wdReplaceAll)\n}\n"}
